$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "catalogBlock_Fish_GenderId" column (C) used to store the text value
# "Hona". The new example is based on a numeric value instead, so update
# C2 and C3 to the number 1.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1

# Update the active selection to reflect the new example cell.
$ws.Range("C4").Select()
